# Applies the "Automatic update of files" edit:
#  1. Updates the "Förändrad" (changed) date in column C for rows 2-11
#     from serial date 45184 to 45186.
#  2. Adds a friendly display-text second argument (the case id from
#     column A) to the HYPERLINK() formulas in columns S, T, V, W, X, Y
#     for the rows that have them (rows 2 and 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "changed" date column (C) for rows 2 through 11 ---
$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = 45186
}

# --- 2. Update the HYPERLINK formulas that include a friendly text arg ---
# Mapping of column letter -> the folder / extension used in the link,
# as found in the existing formulas.
$linkInfo = @(
    @{ Col = "S"; Folder = "artfynd"; Ext = "xlsx" },
    @{ Col = "T"; Folder = "kartor"; Ext = "png" },
    @{ Col = "V"; Folder = "klagomål"; Ext = "docx" },
    @{ Col = "W"; Folder = "klagomålsmail"; Ext = "docx" },
    @{ Col = "X"; Folder = "tillsyn"; Ext = "docx" },
    @{ Col = "Y"; Folder = "tillsynsmail"; Ext = "docx" }
)

$linkRows = @(2, 3)
foreach ($r in $linkRows) {
    $id = $ws.Cells.Item($r, 1).Value2
    foreach ($info in $linkInfo) {
        $col = $info.Col
        $folder = $info.Folder
        $ext = $info.Ext
        $url = "https://klasma.github.io/Logging_TRELLEBORG/" + $folder + "/" + $id + "." + $ext
        $formula = '=HYPERLINK("' + $url + '", "' + $id + '")'
        $ws.Range($col + $r).Formula = $formula
    }
}
